$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, bordered, centered) from H1 into the
# two new header cells I1 and J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2..11 for the new columns I (I0) and J (IF)
$data = @(
    @(1, 4),
    @(1, 6),
    @(3, 7),
    @(1, 6),
    @(1, 5),
    @(2, 5),
    @(3, 3),
    @(8, 9),
    @(7, 7),
    @(1, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
